$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test case row 12 mirrors the formatting of row 11 (same style pattern:
# s=5,3,6,2,5), so copy it down first and then overwrite with the new values.
$ws.Range("A11:E11").Copy($ws.Range("A12:E12"))

$ws.Range("A12").Value = "DRA005"
$ws.Range("B12").Value = "OPQA-4221"
$ws.Range("C12").Value = "Verify that error message "" Incorrect password. Please try again.""should be displayed when user enters incorrect password for existing steam account."
$ws.Range("D12").Value = "Y"

# Row 12 gets an explicit height of 30.
$ws.Rows(12).RowHeight = 30

# Hyperlink B12 -> OPQA-4221 Jira ticket (mirrors the existing B2 hyperlink).
# Adding with TextToDisplay set to the URL makes the OOXML "display" attribute
# match the target URL; the cell's visible text/format is then restored
# (adding a hyperlink otherwise forces the built-in Hyperlink cell style,
# but the target keeps B12 on the same plain style as the rest of the row).
$ws.Hyperlinks.Add($ws.Range("B12"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-4221", "", "", "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-4221")
$ws.Range("B11").Copy($ws.Range("B12"))
$ws.Range("B12").Value = "OPQA-4221"

# Update selection to the newly added row.
$ws.Range("A12:E12").Select()
